# Canopy height and instrument height for France site
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# France site is in row 2 (H1 = "instrument", I1 = "canopy", J1 = "Metadata")
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 0.7

# Add the metadata hyperlink for the France site, matching the pattern
# used by the other sites' "Metadata" column entries.
$ws.Hyperlinks.Add($ws.Range("J2"), "https://meta.icos-cp.eu/resources/stations/ES_FR-EM2") | Out-Null

# Move the active selection as recorded in the saved workbook.
[void]$ws.Range("H3").Select()
